# mascarpone.xlsx -- "Update: fix db and boiling plan"
# Fills in the previously-empty "Выход" (Output/yield, col T) parameter for
# every SKU row, and corrects a batch of "Скорость фасовки" (Filling speed,
# col S) and "Коэффициент" (Coefficient, col U) values for rows 23-53.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("T2").Value = 600
$ws.Range("T3").Value = 600
$ws.Range("T4").Value = 600
$ws.Range("T5").Value = 600
$ws.Range("T6").Value = 600
$ws.Range("T7").Value = 600
$ws.Range("T8").Value = 600
$ws.Range("T9").Value = 600
$ws.Range("T10").Value = 600
$ws.Range("T11").Value = 300
$ws.Range("T12").Value = 300
$ws.Range("T13").Value = 300
$ws.Range("T14").Value = 600
$ws.Range("T15").Value = 600
$ws.Range("T16").Value = 600
$ws.Range("T17").Value = 300
$ws.Range("T18").Value = 600
$ws.Range("T19").Value = 300
$ws.Range("T20").Value = 600
$ws.Range("T21").Value = 600
$ws.Range("T22").Value = 370
$ws.Range("S23").Value = 800
$ws.Range("T23").Value = 370
$ws.Range("S24").Value = 800
$ws.Range("T24").Value = 370
$ws.Range("S25").Value = 800
$ws.Range("T25").Value = 370
$ws.Range("S26").Value = 800
$ws.Range("T26").Value = 370
$ws.Range("S27").Value = 800
$ws.Range("T27").Value = 370
$ws.Range("U27").Value = 1.42
$ws.Range("S28").Value = 800
$ws.Range("T28").Value = 370
$ws.Range("U28").Value = 1.42
$ws.Range("S29").Value = 800
$ws.Range("T29").Value = 370
$ws.Range("S30").Value = 800
$ws.Range("T30").Value = 370
$ws.Range("S31").Value = 800
$ws.Range("T31").Value = 370
$ws.Range("S32").Value = 700
$ws.Range("T32").Value = 370
$ws.Range("S33").Value = 700
$ws.Range("T33").Value = 370
$ws.Range("U33").Value = 1.07
$ws.Range("S34").Value = 700
$ws.Range("T34").Value = 370
$ws.Range("U34").Value = 1.07
$ws.Range("S35").Value = 700
$ws.Range("T35").Value = 370
$ws.Range("U35").Value = 1.07
$ws.Range("S36").Value = 700
$ws.Range("T36").Value = 370
$ws.Range("S37").Value = 800
$ws.Range("T37").Value = 370
$ws.Range("U37").Value = 1.42
$ws.Range("S38").Value = 700
$ws.Range("T38").Value = 370
$ws.Range("S39").Value = 800
$ws.Range("T39").Value = 370
$ws.Range("S40").Value = 800
$ws.Range("T40").Value = 370
$ws.Range("S41").Value = 700
$ws.Range("T41").Value = 370
$ws.Range("U41").Value = 1.07
$ws.Range("S42").Value = 800
$ws.Range("T42").Value = 370
$ws.Range("S43").Value = 700
$ws.Range("T43").Value = 370
$ws.Range("U43").Value = 1.07
$ws.Range("S44").Value = 800
$ws.Range("T44").Value = 370
$ws.Range("U44").Value = 1.42
$ws.Range("S45").Value = 700
$ws.Range("T45").Value = 370
$ws.Range("S46").Value = 700
$ws.Range("T46").Value = 370
$ws.Range("U46").Value = 1.07
$ws.Range("S47").Value = 800
$ws.Range("T47").Value = 370
$ws.Range("S48").Value = 800
$ws.Range("T48").Value = 370
$ws.Range("U48").Value = 1.42
$ws.Range("S49").Value = 700
$ws.Range("T49").Value = 370
$ws.Range("U49").Value = 1.07
$ws.Range("S50").Value = 700
$ws.Range("T50").Value = 370
$ws.Range("U50").Value = 1.07
$ws.Range("T51").Value = 370
$ws.Range("U51").Value = 1.42
$ws.Range("S52").Value = 700
$ws.Range("T52").Value = 370
$ws.Range("U52").Value = 1.42
$ws.Range("S53").Value = 700
$ws.Range("T53").Value = 370
$ws.Range("U53").Value = 1.42

# Match the author's last on-screen selection/scroll position.
$ws.Range("P21").Select()

